# Weekly update: add a new daily price record for
# "Feria Lagunitas de Puerto Montt - Piña" (Fruta / hortaliza, semanal).
#
# The new observation is inserted as row 106 (pushing the existing rows
# 106-138 down to 107-139), matching how the source data is kept in
# (roughly) date order within the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, shifting rows 106:138 down to 107:139
$ws.Rows.Item(106).Insert()

# Populate the new row with the new weekly record
$ws.Range("A106").Value = 4
$ws.Range("B106").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C106").Value = "Los Lagos"
$ws.Range("D106").Value2 = 44463
$ws.Range("E106").Value = 10
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100108
$ws.Range("H106").Value = "Tropicales y subtropicales"
$ws.Range("I106").Value = 100108005
$ws.Range("J106").Value = "Piña"
$ws.Range("K106").Value = "Caramelo"
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 160
$ws.Range("N106").Value = 22000
$ws.Range("O106").Value = 22000
$ws.Range("P106").Value = 22000
$ws.Range("Q106").Value = "$/caja 12 unidades"
$ws.Range("R106").Value = "Ecuador"
$ws.Range("S106").Value = 1833
$ws.Range("T106").Value = 12
